$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells that look numeric as Text so values round-trip exactly
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.688.02"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "2.094.31"
$ws.Range("E3").Value = "  +9.67%  "
$ws.Range("D4").Value = "0.9945"
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").Value = "328.31"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "0.9949"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "0.5170"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("D8").Value = "0.4337"
$ws.Range("E8").Value = "  +6.62%  "
$ws.Range("D9").Value = "0.08841"
$ws.Range("E9").Value = "  +6.30%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.158"
$ws.Range("E10").Value = "  +5.16%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "43.65"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "24.66"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").Value = "2.080.40"
$ws.Range("E13").Value = "  +8.28%  "
$ws.Range("D14").Value = "6.760"
$ws.Range("E14").Value = "  +5.93%  "
$ws.Range("D15").Value = "7.570"
$ws.Range("E15").Value = "  +4.88%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "97.33"
$ws.Range("E16").Value = "  +5.74%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "0.9930"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "0.00001123"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").Value = "0.06589"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "19.13"
$ws.Range("E20").Value = "  +5.02%  "
$ws.Range("D21").Value = "0.9942"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "6.351"
$ws.Range("E22").Value = "  +7.09%  "
$ws.Range("D23").Value = "30.876.85"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("D24").Value = "11.93"
$ws.Range("E24").Value = "  +5.74%  "
$ws.Range("D25").Value = "2.264"
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("D26").Value = "2.316.68"
$ws.Range("E26").Value = "  +8.19%  "
$ws.Range("D27").Value = "22.76"
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.539"
$ws.Range("E28").Value = "  +10.38%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "162.74"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "133.64"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").Value = "1.174"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("D32").Value = "0.1067"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").Value = "6.224"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("D35").Value = "1.473"
$ws.Range("E35").Value = "  +24.22%  "
$ws.Range("D36").Value = "0.02573"
$ws.Range("E36").Value = "  +5.40%  "
$ws.Range("D37").Value = "5.623"
$ws.Range("E37").Value = "  +4.47%  "
$ws.Range("D38").Value = "0.06705"
$ws.Range("E38").Value = "  +5.15%  "
$ws.Range("D39").Value = "9.414"
$ws.Range("E39").Value = "  +9.01%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2238"
$ws.Range("E40").Value = "  +4.30%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "12.40"
$ws.Range("E41").Value = "  +9.69%  "
$ws.Range("D42").Value = "0.6803"
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("D43").Value = "1.248"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "14.09"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "0.9938"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "0.6312"
$ws.Range("E46").Value = "  +4.37%  "
$ws.Range("D47").Value = "2.258"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "3.631"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").Value = "1.267"
$ws.Range("E49").Value = "  +5.32%  "
$ws.Range("D50").Value = "126.97"
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("D51").Value = "82.93"
$ws.Range("E51").Value = "  +5.29%  "
